# Apply the "lowered data scale of early science" edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 (Max) ---
$ws.Range("E3").Value = 65
$ws.Range("G3").Value = 140

# --- Row 4 (Ratio) ---
$ws.Range("E4").Value = 65.001999999999995
$ws.Range("G4").Value = 236.37200000000001

# --- Row 7 (Entry) ---
$ws.Range("E7").Value = 3000
$ws.Range("G7").Value = 4400

# Recalculate so the dependent formulas (F3, H3, F4, H4, E5, G5, F7, H7)
# pick up the new cached values.
$excel.Calculate()

# Move the active selection from I4 to E5, as in the saved workbook view.
$ws.Range("E5").Select()
